$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.378.54"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.880.49"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'0.7233"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").Value = "'243.05"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.08014"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("D9").Value = "'0.3145"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").Value = "'25.03"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "'0.08186"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "1.880.91"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'94.81"
$ws.Range("E13").Value = "  +4.13%  "
$ws.Range("D14").Value = "'5.241"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'0.7135"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "'6.429"
$ws.Range("E16").Value = "  +6.07%  "
$ws.Range("D17").Value = "'0.000008504"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "29.358.38"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'244.66"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").Value = "'13.32"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "2.126.54"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'7.755"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'0.1609"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "'162.84"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'9.054"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "'1.506"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'4.410"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "'4.284"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'1.230"
$ws.Range("E32").Value = "  -5.30%  "
$ws.Range("D33").Value = "'0.05363"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "'1.940"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'0.7650"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").Value = "'1.180"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").Value = "'0.01874"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "1.265.04"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("D40").Value = "'2.759"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "'6.440"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "'113.64"
$ws.Range("E42").Value = "  +4.36%  "
$ws.Range("D43").Value = "'0.9073"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").Value = "'74.48"
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("E45").Value = "  +6.81%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "2.020.79"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.804"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.5202"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'9.501"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "'0.4346"
$ws.Range("E51").Value = "  +0.68%  "
